# Fill in the previously-blank review-count columns on the hotel_info
# sheet (English_Reviews_num, Local_Rank, Total_Reviews_num) for the
# single data row. Orbitz_ReviewURL (J2) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

$reviewCells = $ws.Range("G2:I2")
$reviewCells.NumberFormat = "@"

$ws.Range("G2").Value = "3"
$ws.Range("H2").Value = "423"
$ws.Range("I2").Value = "3"
